$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B183").Value = 7384630
$ws.Range("F183").Value = "Atletico Grau"
$ws.Range("G183").Value = "Unin Comercio"
$ws.Range("H183").Value = 0
$ws.Range("I183").Value = 1
$ws.Range("J183").Value = "A"
$ws.Range("K183").Value = 2.8
$ws.Range("L183").Value = 3.4
$ws.Range("M183").Value = 2.15
$ws.Range("N183").Value = 1.75
$ws.Range("O183").Value = 3.6
$ws.Range("P183").Value = 3.8
$ws.Range("Q183").Value = -0.75
$ws.Range("R183").Value = 2
$ws.Range("S183").Value = 1.8
$ws.Range("T183").Value = 3
$ws.Range("W183").Value = -1
$ws.Range("Y183").Value = 2.8
$ws.Range("Z183").Value = -1
$ws.Range("AA183").Value = 0.8

$ws.Range("B184").Value = 7384627
$ws.Range("F184").Value = "Universitario de Deportes"
$ws.Range("G184").Value = "Sport Huancayo"
$ws.Range("H184").Value = 2
$ws.Range("J184").Value = "H"
$ws.Range("K184").Value = 1.25
$ws.Range("L184").Value = 5
$ws.Range("M184").Value = 12
$ws.Range("N184").Value = 1.181
$ws.Range("O184").Value = 6
$ws.Range("P184").Value = 13
$ws.Range("Q184").Value = -1.75
$ws.Range("R184").Value = 1.8
$ws.Range("S184").Value = 2
$ws.Range("T184").Value = 2.75
$ws.Range("U184").Value = 1.85
$ws.Range("V184").Value = 1.95
$ws.Range("W184").Value = 0.181
$ws.Range("X184").Value = -1
$ws.Range("Z184").Value = 0.4
$ws.Range("AA184").Value = -0.5
$ws.Range("AC184").Value = 0.95

$ws.Range("B185").Value = 7384625
$ws.Range("F185").Value = "AD Tarma"
$ws.Range("G185").Value = "Carlos Manucci"
$ws.Range("I185").Value = 0
$ws.Range("J185").Value = "D"
$ws.Range("K185").Value = 1.5
$ws.Range("L185").Value = 3.75
$ws.Range("M185").Value = 7
$ws.Range("N185").Value = 1.363
$ws.Range("O185").Value = 4.333
$ws.Range("P185").Value = 9.5
$ws.Range("Q185").Value = -1.25
$ws.Range("R185").Value = 1.875
$ws.Range("S185").Value = 1.925
$ws.Range("T185").Value = 2.5
$ws.Range("U185").Value = 1.8
$ws.Range("V185").Value = 2
$ws.Range("X185").Value = 3.333
$ws.Range("Y185").Value = -1
$ws.Range("AA185").Value = 0.925
$ws.Range("AC185").Value = 1

$ws.Range("B186").Value = 7384628
$ws.Range("F186").Value = "Deportivo Binacional"
$ws.Range("G186").Value = "FBC Melgar"
$ws.Range("H186").Value = 1
$ws.Range("I186").Value = 2
$ws.Range("J186").Value = "A"
$ws.Range("K186").Value = 2.75
$ws.Range("L186").Value = 3.3
$ws.Range("M186").Value = 2.375
$ws.Range("N186").Value = 3.3
$ws.Range("O186").Value = 3.6
$ws.Range("P186").Value = 2
$ws.Range("Q186").Value = 0.5
$ws.Range("R186").Value = 1.8
$ws.Range("S186").Value = 2
$ws.Range("T186").Value = 2.75
$ws.Range("U186").Value = 1.975
$ws.Range("V186").Value = 1.875
$ws.Range("W186").Value = -1
$ws.Range("Y186").Value = 1
$ws.Range("Z186").Value = -1
$ws.Range("AA186").Value = 1
$ws.Range("AB186").Value = 0.4875
$ws.Range("AC186").Value = -0.5

$ws.Range("B187").Value = 7384626
$ws.Range("F187").Value = "Sporting Cristal"
$ws.Range("G187").Value = "Alianza Atletico"
$ws.Range("H187").Value = 3
$ws.Range("I187").Value = 0
$ws.Range("J187").Value = "H"
$ws.Range("K187").Value = 1.3
$ws.Range("L187").Value = 5
$ws.Range("M187").Value = 9
$ws.Range("N187").Value = 1.166
$ws.Range("O187").Value = 6.5
$ws.Range("P187").Value = 13
$ws.Range("Q187").Value = -2
$ws.Range("R187").Value = 1.85
$ws.Range("S187").Value = 1.95
$ws.Range("T187").Value = 3.25
$ws.Range("U187").Value = 2
$ws.Range("V187").Value = 1.8
$ws.Range("W187").Value = 0.1659999999999999
$ws.Range("Y187").Value = -1
$ws.Range("Z187").Value = 0.8500000000000001
$ws.Range("AA187").Value = -1
$ws.Range("AB187").Value = -0.5
$ws.Range("AC187").Value = 0.4

$ws.Range("B188").Value = 7384629
$ws.Range("F188").Value = "Deportivo Garcilaso"
$ws.Range("G188").Value = "Alianza Lima"
$ws.Range("H188").Value = 0
$ws.Range("I188").Value = 1
$ws.Range("K188").Value = 2.625
$ws.Range("M188").Value = 2.5
$ws.Range("N188").Value = 2.7
$ws.Range("O188").Value = 3.4
$ws.Range("P188").Value = 2.375
$ws.Range("Q188").Value = 0
$ws.Range("R188").Value = 2.025
$ws.Range("S188").Value = 1.775
$ws.Range("T188").Value = 2.25
$ws.Range("U188").Value = 1.825
$ws.Range("V188").Value = 1.975
$ws.Range("Y188").Value = 1.375
$ws.Range("AA188").Value = 0.7749999999999999
$ws.Range("AB188").Value = -1
$ws.Range("AC188").Value = 0.9750000000000001

$ws.Range("B252").Value = 7883367
$ws.Range("F252").Value = "Sport Boys"
$ws.Range("G252").Value = "Cesar Vallejo"
$ws.Range("I252").Value = 0
$ws.Range("J252").Value = "H"
$ws.Range("K252").Value = 2.2
$ws.Range("L252").Value = 3.3
$ws.Range("M252").Value = 3.1
$ws.Range("N252").Value = 2.4
$ws.Range("O252").Value = 3.25
$ws.Range("P252").Value = 2.625
$ws.Range("Q252").Value = 0
$ws.Range("R252").Value = 1.8
$ws.Range("S252").Value = 2
$ws.Range("T252").Value = 2.5
$ws.Range("U252").Value = 2.025
$ws.Range("W252").Value = 1.4
$ws.Range("X252").Value = -1
$ws.Range("Z252").Value = 0.8
$ws.Range("AA252").Value = -1
$ws.Range("AB252").Value = -1
$ws.Range("AC252").Value = 0.825

$ws.Range("B253").Value = 7882752
$ws.Range("F253").Value = "Sport Huancayo"
$ws.Range("G253").Value = "Union Comercio"
$ws.Range("I253").Value = 2
$ws.Range("J253").Value = "D"
$ws.Range("K253").Value = 1.3
$ws.Range("L253").Value = 4.5
$ws.Range("M253").Value = 10
$ws.Range("N253").Value = 1.3
$ws.Range("O253").Value = 4.333
$ws.Range("P253").Value = 9.5
$ws.Range("Q253").Value = -1.5
$ws.Range("R253").Value = 1.95
$ws.Range("S253").Value = 1.85
$ws.Range("T253").Value = 2.75
$ws.Range("U253").Value = 1.975
$ws.Range("W253").Value = -1
$ws.Range("X253").Value = 3.333
$ws.Range("Z253").Value = -1
$ws.Range("AA253").Value = 0.8500000000000001
$ws.Range("AB253").Value = 0.9750000000000001
$ws.Range("AC253").Value = -1

$ws.Range("B265").Value = 7971191
$ws.Range("E265").Value = 45381.91666666666
$ws.Range("F265").Value = "Cesar Vallejo"
$ws.Range("G265").Value = "Universitario de Deportes"
$ws.Range("K265").Value = 3.75
$ws.Range("L265").Value = 3.4
$ws.Range("M265").Value = 1.95
$ws.Range("N265").Value = 2.875
$ws.Range("O265").Value = 3.2
$ws.Range("P265").Value = 2.45
$ws.Range("Q265").Value = 0
$ws.Range("R265").Value = 2.1
$ws.Range("S265").Value = 1.775
$ws.Range("T265").Value = 2.25
$ws.Range("U265").Value = 1.975
$ws.Range("V265").Value = 1.875

$ws.Range("B266").Value = 7971192
$ws.Range("E266").Value = 45382.625
$ws.Range("F266").Value = "Sport Huancayo"
$ws.Range("G266").Value = "AD Tarma"
$ws.Range("K266").Value = 1.8
$ws.Range("L266").Value = 3.75
$ws.Range("M266").Value = 4
$ws.Range("N266").Value = 2.05
$ws.Range("O266").Value = 3.6
$ws.Range("P266").Value = 3.4
$ws.Range("Q266").Value = -0.25
$ws.Range("R266").Value = 1.825
$ws.Range("S266").Value = 2.025
$ws.Range("U266").Value = 1.975
$ws.Range("V266").Value = 1.875

$ws.Range("B267").Value = 7971193
$ws.Range("E267").Value = 45382.70833333334
$ws.Range("F267").Value = "Sport Boys"
$ws.Range("G267").Value = "Alianza Atletico"
$ws.Range("K267").Value = 1.615
$ws.Range("L267").Value = 4
$ws.Range("M267").Value = 5
$ws.Range("N267").Value = 1.909
$ws.Range("O267").Value = 3.8
$ws.Range("P267").Value = 3.75
$ws.Range("Q267").Value = -0.5
$ws.Range("R267").Value = 1.925
$ws.Range("S267").Value = 1.925
$ws.Range("U267").Value = 2
$ws.Range("V267").Value = 1.85

$ws.Range("B268").Value = 7971194
$ws.Range("E268").Value = 45382.72916666666
$ws.Range("F268").Value = "UTC Cajamarca"
$ws.Range("G268").Value = "Sporting Cristal"
$ws.Range("K268").Value = 1.8
$ws.Range("L268").Value = 3.6
$ws.Range("M268").Value = 1.8
$ws.Range("N268").Value = 4.2
$ws.Range("R268").Value = 1.925
$ws.Range("S268").Value = 1.925
$ws.Range("T268").Value = 2.5
$ws.Range("U268").Value = 1.925
$ws.Range("V268").Value = 1.925

Write-Output "Done applying Peru Liga 1 odds/id updates for rows 183-188, 252-253, 265-268."
